$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# New translation table content for the "Translation" sheet (rows 4-28).
# The GPS/overview screen texts (Draw, GpsData) were promoted next to the
# header, and new rows were added for the Debug/Application screen split
# plus new GpsData fields (Date, Read, Write) and two placeholder rows.
$data = @(
  @("SingleUseId1", "Large", "Left", "LTR", "pathTracker"),
  @("SingleUseId17", "Default", "Center", "LTR", "Draw"),
  @("SingleUseId18", "Default", "Center", "LTR", "GpsData"),
  @("SingleUseId25", "Default", "Center", "LTR", "Debug"),
  @("SingleUseId26", "Default", "Center", "LTR", "Application"),
  @("SingleUseId2", "Small", "Left", "LTR", "Time: <value>"),
  @("SingleUseId3", "Small", "Left", "LTR", "00:00:00"),
  @("SingleUseId4", "Small", "Left", "LTR", "Lat: <value> <value>"),
  @("SingleUseId8", "Small", "Left", "LTR", "---"),
  @("SingleUseId9", "Small", "Left", "LTR", "Lon: <value> <value>"),
  @("SingleUseId10", "Small", "Left", "LTR", "---"),
  @("SingleUseId11", "Small", "Left", "LTR", "Alti: <value> m"),
  @("SingleUseId12", "Small", "Left", "LTR", "---"),
  @("SingleUseId13", "Small", "Left", "LTR", "Fix: <value>"),
  @("SingleUseId14", "Small", "Left", "LTR", "---"),
  @("SingleUseId15", "Small", "Left", "LTR", "Sat: <value>"),
  @("SingleUseId16", "Small", "Left", "LTR", "---"),
  @("SingleUseId19", "Small", "Left", "LTR", "Date: <value>"),
  @("SingleUseId20", "Small", "Left", "LTR", "01.01.2020"),
  @("SingleUseId21", "Small", "Left", "LTR", "Read: <value>"),
  @("SingleUseId22", "Small", "Left", "LTR", "---"),
  @("SingleUseId23", "Small", "Left", "LTR", "Write: <value>"),
  @("SingleUseId24", "Small", "Left", "LTR", "---"),
  @("SingleUseId28", "Default", "Center", "LTR", "Application"),
  @("SingleUseId29", "Default", "Center", "LTR", "New Text")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = 4 + $i
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, 2 + $j).Value = $row[$j]
    }
}
